$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("is_base") values change from numeric 1/0 to text "TRUE"/"FALSE"
$ws.Range("D2").Formula = "'TRUE"
$ws.Range("D3").Formula = "'TRUE"
$ws.Range("D4").Formula = "'TRUE"
$ws.Range("D5").Formula = "'TRUE"
$ws.Range("D6").Formula = "'FALSE"
$ws.Range("D7").Formula = "'FALSE"

# Remove the quote-prefix cell formatting that entering text this way applies,
# so the cells keep their original (unstyled) appearance.
$ws.Range("D2:D7").Style = "Normal"
